$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells stay formatted as Text so that
# numeric-looking strings (e.g. "1.00", "0.260", thousand-dot
# separated prices) are preserved exactly as written, matching
# the original inlineStr cell contents.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.844.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.093.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.087.22"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.95"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.593.36"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.858.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.094.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.24"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "502.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.306.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.260"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0543"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -14.75%  "
